$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7764291763305664
$ws.Range("B1").Value = 3.034293413162231
$ws.Range("C1").Value = 3.996899366378784
$ws.Range("D1").Value = 0.9416834115982056
$ws.Range("E1").Value = 0.8759952187538147
